# Updated the single-record upload fixture to match the new data-upload
# template: the "validation" column is dropped, "notes" slides into its
# place (AD), and two new trailing columns are appended - "hatchery" (AE)
# and "agency_stock_id" (AF) - with sample values "CFCW" and "X9999".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row (row 1): AD was "validation", now becomes "notes" (shifted
# up from the old AE), and two brand-new headers are appended after it.
$ws.Range("AD1").Value = "notes"
$ws.Range("AE1").Value = "hatchery"
$ws.Range("AF1").Value = "agency_stock_id"

# Data row (row 2): the old AE2 value ("FIS ID = 73699") now lives under
# the "notes" header at AD2; the new hatchery/agency_stock_id columns get
# their sample values.
$ws.Range("AD2").Value = "FIS ID = 73699"
$ws.Range("AE2").Value = "CFCW"
$ws.Range("AF2").Value = "X9999"

# Reflect the new used range in the view: scrolled right to show the new
# columns, with AF5 as the active selection.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 26
[void]$ws.Range("AF5").Select()
